$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$p1 = $d.Paragraphs.Item($count - 1)
$p2 = $d.Paragraphs.Item($count)

# Set the later paragraph first, then the earlier one, so that neither
# assignment needs an embedded paragraph mark and no index shifting occurs.
$p2.Range.Text = "We create a new function add to the posts controller"
$p1.Range.Text = "We want to be able to add a post"
